$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 18 of user history data
$ws.Range("A18").Value = "nam354"

$ws.Range("B18").Value = 200
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0

# Copy the number-format/alignment style used by the row above (B17:E17) down to B18:E18
$ws.Range("B17:E17").Copy()
$ws.Range("B18:E18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F18").Value = ";1"
$ws.Range("G18").Value = ";0"
$ws.Range("H18").Value = ";+100"
